$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9142154455184937
$ws.Range("B1").Value = 1.721049785614014
$ws.Range("C1").Value = 4.042654037475586
$ws.Range("D1").Value = 3.733898639678955
$ws.Range("E1").Value = 0.692113995552063
